# Generate Report for Handoff
#
# Refreshes the localization handoff report: for every file row whose
# Status is "Ready for handoff" (rows 7, 8, 9, 10, 12, 13 on each language
# table; row 11 already has a newer handback and is left untouched), stamp
# the latest handoff generation time and flag the row's Priority as "ht".

$wb = $excel.ActiveWorkbook

$rows = 7, 8, 9, 10, 12, 13

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-13 10:25:33"
}

# --- zh-cn sheet: "Latest Handoff Datetime" (H) and "Priority" (E) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-08-13 10:25:25"
    $wsZhCn.Range("E$r").Value = "ht"
}

# --- de-de sheet: "Latest Handoff Datetime" (H) and "Priority" (E) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("H$r").Value = "2016-08-13 10:25:33"
    $wsDeDe.Range("E$r").Value = "ht"
}
